$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclaimer text: change the model-holdings date
# from 2021-03-29 to 2021-03-30.
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."
$ws.Rows("18").AutoFit()

# Update the Weight (D) and Percent Change (E) columns for rows 2-15.
$ws.Range("D2").Value = 0.05580023786378768
$ws.Range("E2").Value = -0.000291970802919761

$ws.Range("D3").Value = 0.02343134805771727
$ws.Range("E3").Value = -0.001216791725816346

$ws.Range("D4").Value = 0.03153120740310309
$ws.Range("E4").Value = 0.00406897888006208

$ws.Range("D5").Value = 0.03232205019034534
$ws.Range("E5").Value = -0.0103750997605746

$ws.Range("D6").Value = 0.03387115168432543
$ws.Range("E6").Value = 0.003367003367003463

$ws.Range("D7").Value = 0.01894560630827799
$ws.Range("E7").Value = 0.004765487835465487

$ws.Range("D8").Value = 0.004563807775647745
$ws.Range("E8").Value = 0.04819277108433728

$ws.Range("D9").Value = 0.006660729123108088
$ws.Range("E9").Value = 0.0008153281695881809

$ws.Range("D10").Value = 0.06924117107185333
$ws.Range("E10").Value = -0.01647058823529424

$ws.Range("D11").Value = 0.06936336137374484
$ws.Range("E11").Value = -0.01702877275396375

$ws.Range("D12").Value = 0.1471714302782137
$ws.Range("E12").Value = 0.005239852398523981

$ws.Range("D13").Value = 0.3918704076811539
$ws.Range("E13").Value = -0.0008834702712252618

$ws.Range("D14").Value = 0.1152274911887216
$ws.Range("E14").Value = 0.0002474328838302231

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -0.001690299176165699

# Restore sheet protection (the worksheet was protected before this edit).
$ws.Protect()
